$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("data")

# ---- 1. Update query timestamps in column F of the "data" sheet ----
$newTimes = @(
    "2021-10-05 14:33:22.335539",
    "2021-10-05 14:33:22.335549",
    "2021-10-05 14:33:22.335553",
    "2021-10-05 14:33:22.335555",
    "2021-10-05 14:33:22.335558",
    "2021-10-05 14:33:22.335561",
    "2021-10-05 14:33:22.335564",
    "2021-10-05 14:33:22.335566",
    "2021-10-05 14:33:22.335570",
    "2021-10-05 14:33:22.335572",
    "2021-10-05 14:33:22.335574",
    "2021-10-05 14:33:22.335577",
    "2021-10-05 14:33:22.335580",
    "2021-10-05 14:33:22.335582",
    "2021-10-05 14:33:22.335585",
    "2021-10-05 14:33:22.335587",
    "2021-10-05 14:33:22.335590",
    "2021-10-05 14:33:22.335593",
    "2021-10-05 14:33:22.335596",
    "2021-10-05 14:33:22.335599",
    "2021-10-05 14:33:22.335601",
    "2021-10-05 14:33:22.335604",
    "2021-10-05 14:33:22.335606",
    "2021-10-05 14:33:22.335609",
    "2021-10-05 14:33:22.335612",
    "2021-10-05 14:33:22.335614",
    "2021-10-05 14:33:22.335617",
    "2021-10-05 14:33:22.335620",
    "2021-10-05 14:33:22.335622",
    "2021-10-05 14:33:22.335625",
    "2021-10-05 14:33:22.335627",
    "2021-10-05 14:33:22.335630",
    "2021-10-05 14:33:22.335633",
    "2021-10-05 14:33:22.335636",
    "2021-10-05 14:33:22.335638",
    "2021-10-05 14:33:22.335641",
    "2021-10-05 14:33:22.335643",
    "2021-10-05 14:33:22.335646",
    "2021-10-05 14:33:22.335648",
    "2021-10-05 14:33:22.335651",
    "2021-10-05 14:33:22.335654",
    "2021-10-05 14:33:22.335657",
    "2021-10-05 14:33:22.335660",
    "2021-10-05 14:33:22.335662",
    "2021-10-05 14:33:22.335665",
    "2021-10-05 14:33:22.335668",
    "2021-10-05 14:33:22.335670",
    "2021-10-05 14:33:22.335672",
    "2021-10-05 14:33:22.335675",
    "2021-10-05 14:33:22.335678",
    "2021-10-05 14:33:22.335680",
    "2021-10-05 14:33:22.335683",
    "2021-10-05 14:33:22.335686",
    "2021-10-05 14:33:22.335689",
    "2021-10-05 14:33:22.335691",
    "2021-10-05 14:33:22.335694",
    "2021-10-05 14:33:22.335697",
    "2021-10-05 14:33:22.335700",
    "2021-10-05 14:33:22.335702",
    "2021-10-05 14:33:22.335705",
    "2021-10-05 14:33:22.335708",
    "2021-10-05 14:33:22.335710",
    "2021-10-05 14:33:22.335713",
    "2021-10-05 14:33:22.335715",
    "2021-10-05 14:33:22.335719",
    "2021-10-05 14:33:22.335722",
    "2021-10-05 14:33:22.335725",
    "2021-10-05 14:33:22.335727",
    "2021-10-05 14:33:22.335730",
    "2021-10-05 14:33:22.335733",
    "2021-10-05 14:33:22.335735",
    "2021-10-05 14:33:22.335738",
    "2021-10-05 14:33:22.335740",
    "2021-10-05 14:33:22.335743",
    "2021-10-05 14:33:22.335745",
    "2021-10-05 14:33:22.335748",
    "2021-10-05 14:33:22.335752",
    "2021-10-05 14:33:22.335755",
    "2021-10-05 14:33:22.335758",
    "2021-10-05 14:33:22.335761",
    "2021-10-05 14:33:22.335763",
    "2021-10-05 14:33:22.335766",
    "2021-10-05 14:33:22.335768",
    "2021-10-05 14:33:22.335771",
    "2021-10-05 14:33:22.335774",
    "2021-10-05 14:33:22.335776",
    "2021-10-05 14:33:22.335779",
    "2021-10-05 14:33:22.335781",
    "2021-10-05 14:33:22.335784",
    "2021-10-05 14:33:22.335786",
    "2021-10-05 14:33:22.335789",
    "2021-10-05 14:33:22.335791",
    "2021-10-05 14:33:22.335795",
    "2021-10-05 14:33:22.335798",
    "2021-10-05 14:33:22.335801",
    "2021-10-05 14:33:22.335803",
    "2021-10-05 14:33:22.335806",
    "2021-10-05 14:33:22.335809",
    "2021-10-05 14:33:22.335811",
    "2021-10-05 14:33:22.335814",
    "2021-10-05 14:33:22.335817",
    "2021-10-05 14:33:22.335819",
    "2021-10-05 14:33:22.335822",
    "2021-10-05 14:33:22.335824",
    "2021-10-05 14:33:22.335827",
    "2021-10-05 14:33:22.335830",
    "2021-10-05 14:33:22.335833",
    "2021-10-05 14:33:22.335835",
    "2021-10-05 14:33:22.335840",
    "2021-10-05 14:33:22.335843",
    "2021-10-05 14:33:22.335846",
    "2021-10-05 14:33:22.335848",
    "2021-10-05 14:33:22.335851",
    "2021-10-05 14:33:22.335853",
    "2021-10-05 14:33:22.335856",
    "2021-10-05 14:33:22.335858",
    "2021-10-05 14:33:22.335861",
    "2021-10-05 14:33:22.335864",
    "2021-10-05 14:33:22.335866",
    "2021-10-05 14:33:22.335869",
    "2021-10-05 14:33:22.335872",
    "2021-10-05 14:33:22.335874",
    "2021-10-05 14:33:22.335879",
    "2021-10-05 14:33:22.335896",
    "2021-10-05 14:33:22.335905",
    "2021-10-05 14:33:22.335909",
    "2021-10-05 14:33:22.335911",
    "2021-10-05 14:33:22.335914",
    "2021-10-05 14:33:22.335920",
    "2021-10-05 14:33:22.335923",
    "2021-10-05 14:33:22.335926",
    "2021-10-05 14:33:22.335928",
    "2021-10-05 14:33:22.335932",
    "2021-10-05 14:33:22.335934",
    "2021-10-05 14:33:22.335937",
    "2021-10-05 14:33:22.335939",
    "2021-10-05 14:33:22.335942",
    "2021-10-05 14:33:22.335944",
    "2021-10-05 14:33:22.335947",
    "2021-10-05 14:33:22.335949",
    "2021-10-05 14:33:22.335952",
    "2021-10-05 14:33:22.335955",
    "2021-10-05 14:33:22.335957",
    "2021-10-05 14:33:22.335960",
    "2021-10-05 14:33:22.335963",
    "2021-10-05 14:33:22.335965",
    "2021-10-05 14:33:22.335968",
    "2021-10-05 14:33:22.335970",
    "2021-10-05 14:33:22.335973",
    "2021-10-05 14:33:22.335976",
    "2021-10-05 14:33:22.335979",
    "2021-10-05 14:33:22.335981",
    "2021-10-05 14:33:22.335984",
    "2021-10-05 14:33:22.335986",
    "2021-10-05 14:33:22.335989",
    "2021-10-05 14:33:22.335992",
    "2021-10-05 14:33:22.335995",
    "2021-10-05 14:33:22.335997",
    "2021-10-05 14:33:22.336000",
    "2021-10-05 14:33:22.336002",
    "2021-10-05 14:33:22.336005",
    "2021-10-05 14:33:22.336008",
    "2021-10-05 14:33:22.336010",
    "2021-10-05 14:33:22.336013",
    "2021-10-05 14:33:22.336015",
    "2021-10-05 14:33:22.336018",
    "2021-10-05 14:33:22.336020",
    "2021-10-05 14:33:22.336022",
    "2021-10-05 14:33:22.336025",
    "2021-10-05 14:33:22.336028",
    "2021-10-05 14:33:22.336030",
    "2021-10-05 14:33:22.336032",
    "2021-10-05 14:33:22.336037",
    "2021-10-05 14:33:22.336040",
    "2021-10-05 14:33:22.336042",
    "2021-10-05 14:33:22.336045",
    "2021-10-05 14:33:22.336047",
    "2021-10-05 14:33:22.336050",
    "2021-10-05 14:33:22.336052",
    "2021-10-05 14:33:22.336055",
    "2021-10-05 14:33:22.336057",
    "2021-10-05 14:33:22.336060",
    "2021-10-05 14:33:22.336062",
    "2021-10-05 14:33:22.336065",
    "2021-10-05 14:33:22.336067",
    "2021-10-05 14:33:22.336070",
    "2021-10-05 14:33:22.336073",
    "2021-10-05 14:33:22.336075",
    "2021-10-05 14:33:22.336078",
    "2021-10-05 14:33:22.336080",
    "2021-10-05 14:33:22.336083",
    "2021-10-05 14:33:22.336086",
    "2021-10-05 14:33:22.336088",
    "2021-10-05 14:33:22.336091",
    "2021-10-05 14:33:22.336093",
    "2021-10-05 14:33:22.336096",
    "2021-10-05 14:33:22.336098",
    "2021-10-05 14:33:22.336101",
    "2021-10-05 14:33:22.336104",
    "2021-10-05 14:33:22.336106",
    "2021-10-05 14:33:22.336109",
    "2021-10-05 14:33:22.336112",
    "2021-10-05 14:33:22.336114",
    "2021-10-05 14:33:22.336117",
    "2021-10-05 14:33:22.336119",
    "2021-10-05 14:33:22.336122",
    "2021-10-05 14:33:22.336125",
    "2021-10-05 14:33:22.336127",
    "2021-10-05 14:33:22.336130",
    "2021-10-05 14:33:22.336133",
    "2021-10-05 14:33:22.336136",
    "2021-10-05 14:33:22.336139",
    "2021-10-05 14:33:22.336141",
    "2021-10-05 14:33:22.336143",
    "2021-10-05 14:33:22.336146",
    "2021-10-05 14:33:22.336150",
    "2021-10-05 14:33:22.336152",
    "2021-10-05 14:33:22.336155",
    "2021-10-05 14:33:22.336157",
    "2021-10-05 14:33:22.336160",
    "2021-10-05 14:33:22.336162",
    "2021-10-05 14:33:22.336164",
    "2021-10-05 14:33:22.336167"
)

for ($i = 0; $i -lt $newTimes.Count; $i++) {
    $row = $i + 2
    $ds.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ---- 2. Add a new "metadata" worksheet after "data" ----
$newSheet = $wb.Worksheets.Add($null, $ds)
$newSheet.Name = "metadata"

# Copy the header cell style (bold, bordered, centered) from the data sheet
# so the new sheet reuses the same existing style instead of creating a new one.
$ds.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Cardiomyopathy_Paediatric"
$newSheet.Range("C2").Value = 3270
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.108"
$newSheet.Range("E2").Value = "2021-09-22T05:52:47.810993Z"
$newSheet.Range("F2").Value = "2021-10-05 14:33:22.331922"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3270/?format=json"
